$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "256.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.40%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.94%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.639"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-10.70%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05897"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.64%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.645"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.95%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8676"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.31%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9459"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.07%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1402"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.63%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.03742"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "8.17%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07076"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.66%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03205"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.19%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09252"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.54%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001545"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.00%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006038"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-94.27%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006010"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.48%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.513"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.44%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.191"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.31%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.223"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.09%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3106"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.24%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.98%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.846"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.12%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04242"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.35%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.37%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001220"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.39%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.19%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001199"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.06%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001501"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2.35%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03814"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.08%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006237"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "10.37%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.22%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002199"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.19%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01144"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "16.25%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005499"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.55%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.06%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06019"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-33.14%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002279"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "7.07%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.06%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
